$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 45.90594266666667
$ws.Range("N2").Value = 137.717828
$ws.Range("O2").Value = 0.3954672001633582
$ws.Range("P2").Value = 0.3954672001633583
$ws.Range("Q2").Value = 1640.625660280524
$ws.Range("R2").Value = 14765.63094252471
$ws.Range("S2").Value = 0.00771054030422371
$ws.Range("T2").Value = 0.007710540304223711

# Row 3
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.3484294080560655
$ws.Range("P3").Value = 0.3484294080560656
$ws.Range("Q3").Value = 1445.48581378431
$ws.Range("R3").Value = 13009.37232405879
$ws.Range("S3").Value = 0.006793430638200438
$ws.Range("T3").Value = 0.006793430638200438

# Row 4
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 12.761795
$ws.Range("N4").Value = 38.28538500000001
$ws.Range("O4").Value = 0.1099393900775594
$ws.Range("P4").Value = 0.1099393900775594
$ws.Range("Q4").Value = 456.0918942514767
$ws.Range("R4").Value = 4104.827048263291
$ws.Range("S4").Value = 0.002143520620331174
$ws.Range("T4").Value = 0.002143520620331174

# Row 5
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 16.966758
$ws.Range("N5").Value = 50.900274
$ws.Range("O5").Value = 0.1461640017030168
$ws.Range("P5").Value = 0.1461640017030168
$ws.Range("Q5").Value = 606.3724417706441
$ws.Range("R5").Value = 5457.351975935796
$ws.Range("S5").Value = 0.002849802526460337
$ws.Range("T5").Value = 0.002849802526460337

# Row 6
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 45.90594266666667
$ws.Range("N6").Value = 137.717828
$ws.Range("O6").Value = 0.3954672001633582
$ws.Range("P6").Value = 0.3954672001633583
$ws.Range("Q6").Value = 77548.41805925308
$ws.Range("R6").Value = 697935.7625332777
$ws.Range("S6").Value = 0.364458643705732
$ws.Range("T6").Value = 0.3644586437057321

# Row 7
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.3484294080560655
$ws.Range("P7").Value = 0.3484294080560656
$ws.Range("Q7").Value = 68324.62815856401
$ws.Range("R7").Value = 614921.6534270761
$ws.Range("S7").Value = 0.3211090817009575
$ws.Range("T7").Value = 0.3211090817009575

# Row 8
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 12.761795
$ws.Range("N8").Value = 38.28538500000001
$ws.Range("O8").Value = 0.1099393900775594
$ws.Range("P8").Value = 0.1099393900775594
$ws.Range("Q8").Value = 21558.36382737213
$ws.Range("R8").Value = 194025.2744463492
$ws.Range("S8").Value = 0.1013190499261416
$ws.Range("T8").Value = 0.1013190499261416

# Row 9
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 16.966758
$ws.Range("N9").Value = 50.900274
$ws.Range("O9").Value = 0.1461640017030168
$ws.Range("P9").Value = 0.1461640017030168
$ws.Range("Q9").Value = 28661.76285820112
$ws.Range("R9").Value = 257955.86572381
$ws.Range("S9").Value = 0.1347032922004124
$ws.Range("T9").Value = 0.1347032922004124

# Row 10
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 45.90594266666667
$ws.Range("N10").Value = 137.717828
$ws.Range("O10").Value = 0.3954672001633582
$ws.Range("P10").Value = 0.3954672001633583
$ws.Range("Q10").Value = 4298.703763235628
$ws.Range("R10").Value = 38688.33386912065
$ws.Range("S10").Value = 0.02020285883903527
$ws.Range("T10").Value = 0.02020285883903527

# Row 11
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.3484294080560655
$ws.Range("P11").Value = 0.3484294080560656
$ws.Range("Q11").Value = 3787.405901207148
$ws.Range("R11").Value = 34086.65311086433
$ws.Range("S11").Value = 0.01779988363995181
$ws.Range("T11").Value = 0.01779988363995182

# Row 12
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 12.761795
$ws.Range("N12").Value = 38.28538500000001
$ws.Range("O12").Value = 0.1099393900775594
$ws.Range("P12").Value = 0.1099393900775594
$ws.Range("Q12").Value = 1195.034302867635
$ws.Range("R12").Value = 10755.30872580872
$ws.Range("S12").Value = 0.005616369644989741
$ws.Range("T12").Value = 0.005616369644989743

# Row 13
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 16.966758
$ws.Range("N13").Value = 50.900274
$ws.Range("O13").Value = 0.1461640017030168
$ws.Range("P13").Value = 0.1461640017030168
$ws.Range("Q13").Value = 1588.793568495174
$ws.Range("R13").Value = 14299.14211645657
$ws.Range("S13").Value = 0.007466942119434362
$ws.Range("T13").Value = 0.007466942119434362

# Row 14
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 45.90594266666667
$ws.Range("N14").Value = 137.717828
$ws.Range("O14").Value = 0.3954672001633582
$ws.Range("P14").Value = 0.3954672001633583
$ws.Range("Q14").Value = 658.5782983034566
$ws.Range("R14").Value = 5927.204684731108
$ws.Range("S14").Value = 0.003095157314367253
$ws.Range("T14").Value = 0.003095157314367253

# Row 15
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.3484294080560655
$ws.Range("P15").Value = 0.3484294080560656
$ws.Range("Q15").Value = 580.2454578828698
$ws.Range("R15").Value = 5222.209120945828
$ws.Range("S15").Value = 0.002727012076955822
$ws.Range("T15").Value = 0.002727012076955823

# Row 16
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 12.761795
$ws.Range("N16").Value = 38.28538500000001
$ws.Range("O16").Value = 0.1099393900775594
$ws.Range("P16").Value = 0.1099393900775594
$ws.Range("Q16").Value = 183.0839483119984
$ws.Range("R16").Value = 1647.755534807985
$ws.Range("S16").Value = 0.0008604498860969279
$ws.Range("T16").Value = 0.0008604498860969281

# Row 17
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 16.966758
$ws.Range("N17").Value = 50.900274
$ws.Range("O17").Value = 0.1461640017030168
$ws.Range("P17").Value = 0.1461640017030168
$ws.Range("Q17").Value = 243.409414168946
$ws.Range("R17").Value = 2190.684727520514
$ws.Range("S17").Value = 0.001143964856709745
$ws.Range("T17").Value = 0.001143964856709745
